$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "bleu" -> "noir" (status color label)
$ws.Range("B5").Value = "noir"
$ws.Range("B6").Value = "noir"
$ws.Range("B7").Value = "noir"

# "résultat et / ou publication posté..." -> "résultat postés ou publiés..."
$ws.Range("C2").Value = "résultat postés ou publiés"
$ws.Range("C3").Value = "résultat postés ou publiés dans les 12 mois"
$ws.Range("C4").Value = "résultat postés ou publiés dans les 36 mois"

# "pas de résultat ni de publication" -> "pas de résultat postés ni publiés"
$ws.Range("C5").Value = "pas de résultat postés ni publiés"
$ws.Range("C6").Value = "pas de résultat postés ni publiés"
$ws.Range("C7").Value = "pas de résultat postés ni publiés"
